$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. Insert a new "FortMill" worksheet between "Denver" and "Huntersville".
#    Easiest/most faithful way is to copy the existing "Denver" sheet (which
#    already has the correct headers/category rows and no location-specific
#    data) and rename the copy.
# ---------------------------------------------------------------------------
$denver = $wb.Worksheets.Item("Denver")
$denver.Copy([System.Reflection.Missing]::Value, $denver)
$fortMill = $wb.Worksheets.Item("Denver (2)")
$fortMill.Name = "FortMill"

# ---------------------------------------------------------------------------
# 2. Recreate the "Summary" sheet (copy it, delete the original, rename the
#    copy back to "Summary", move it back to the front). This mirrors what
#    happened in the authentic edit (Summary's internal sheetId moved from
#    20 to 22) while keeping its contents/position identical.
# ---------------------------------------------------------------------------
$oldSummary = $wb.Worksheets.Item("Summary")
$oldSummary.Copy($oldSummary)
$newSummary = $wb.Worksheets.Item("Summary (2)")
$wb.Worksheets.Item("Summary").Delete()
$newSummary.Name = "Summary"
$newSummary.Move($wb.Worksheets.Item("Ballantyne"))

$summary = $wb.Worksheets.Item("Summary")

# ---------------------------------------------------------------------------
# 3. Roll every month header (Jan..Dec) from 24 to 25 on every sheet
#    (Summary + each location sheet, including the new FortMill one).
# ---------------------------------------------------------------------------
$months24 = @("Jan24","Feb24","Mar24","Apr24","May24","Jun24","Jul24","Aug24","Sep24","Oct24","Nov24","Dec24")
$months25 = @("Jan25","Feb25","Mar25","Apr25","May25","Jun25","Jul25","Aug25","Sep25","Oct25","Nov25","Dec25")
$cols     = @("C","D","E","F","G","H","I","J","K","L","M","N")

foreach ($ws in $wb.Worksheets) {
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $cell = $ws.Range($cols[$i] + "1")
        if ($cell.Value2 -eq $months24[$i]) {
            $cell.Value = $months25[$i]
        }
    }
}

# ---------------------------------------------------------------------------
# 4. Extend the Summary roll-up formulas (rows 2-10, cols C-N) so they also
#    add in the new FortMill location.
# ---------------------------------------------------------------------------
for ($r = 2; $r -le 10; $r++) {
    foreach ($col in $cols) {
        $addr = $col + $r
        $cell = $summary.Range($addr)
        $f = $cell.Formula
        if ($f -ne $null -and $f -ne "") {
            $cell.Formula = $f + "+FortMill!" + $addr
        }
    }
}
